# "git commit calificacion hasta p23"
#
# Adds a new trailing column ("l1") to Sheet1, right after the existing
# "c1" column (P), mirroring the sibling "c1"/"nota_iniciativa" columns:
#   - R1 header cell gets the new label "l1" with the same header
#     formatting (style) as the other header cells.
#   - R2:R78 get the same default numeric value (0) used throughout the
#     "c1" (P) and "nota_iniciativa" (Q) columns for every submission row.
#
# This grows the sheet's used range from A1:Q78 to A1:R78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 78
$newCol = 18   # column R

# Give the new header cell (R1) the same look as the rest of row 1's
# header cells (bold, centered, bordered) by copying Q1's formatting,
# then set its text to the new header label.
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "l1"

# Populate the new column's data rows with 0, the same default used by
# the neighboring "c1"/"nota_iniciativa" columns on every row.
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $newCol).Value = 0
}

Write-Host "Added column R (l1) for rows 1-$lastRow"
